$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CertRenewalWindowOpens (H) and CertRenewalDeadline (I) for the rows
# whose renewal window shifted forward (new H = old I, new I = old I + 183).

$ws.Range("H2").Value = 44372
$ws.Range("I2").Value = 44555

$ws.Range("H18").Value = 44363
$ws.Range("I18").Value = 44546

$ws.Range("H66").Value = 44349
$ws.Range("I66").Value = 44532

$ws.Range("H82").Value = 44355
$ws.Range("I82").Value = 44538

$ws.Range("H98").Value = 44374
$ws.Range("I98").Value = 44557

$ws.Range("H114").Value = 44349
$ws.Range("I114").Value = 44532

$ws.Range("H130").Value = 44374
$ws.Range("I130").Value = 44557

$ws.Range("H146").Value = 44352
$ws.Range("I146").Value = 44535
